# Update "想去人数" (interested-count) figures across the four sheets to
# reflect newly scraped attendance numbers (gh-pages data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 54939
$ws.Range("F4").Value = 3109
$ws.Range("F6").Value = 1171
$ws.Range("F8").Value = 859
$ws.Range("F10").Value = 1095
$ws.Range("F11").Value = 1357
$ws.Range("F12").Value = 118
$ws.Range("F14").Value = 216
$ws.Range("F15").Value = 395
$ws.Range("F16").Value = 52
$ws.Range("F21").Value = 5374
$ws.Range("F23").Value = 5267
$ws.Range("F24").Value = 9262
$ws.Range("F26").Value = 160
$ws.Range("F27").Value = 148
$ws.Range("F28").Value = 239
$ws.Range("F30").Value = 139
$ws.Range("F31").Value = 105
$ws.Range("F32").Value = 4269
$ws.Range("F33").Value = 282

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G3").Value = 249
$ws.Range("F10").Value = 1148

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 591

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 591
$ws.Range("F5").Value = 3109
$ws.Range("G7").Value = 249
$ws.Range("F8").Value = 1171
$ws.Range("F11").Value = 859
$ws.Range("F13").Value = 1095
$ws.Range("F15").Value = 1357
$ws.Range("F17").Value = 118
$ws.Range("F18").Value = 216
$ws.Range("F20").Value = 395
$ws.Range("F21").Value = 52
$ws.Range("F26").Value = 5374
$ws.Range("F28").Value = 5267
$ws.Range("F29").Value = 9262
$ws.Range("F32").Value = 160
$ws.Range("F33").Value = 148
$ws.Range("F34").Value = 239
$ws.Range("F39").Value = 139
$ws.Range("F40").Value = 105
$ws.Range("F41").Value = 4269
$ws.Range("F47").Value = 282

